$wb = $excel.ActiveWorkbook

# 1. Metadata: update the Date value (row 8, column B)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(8, 2).Value = "2023-09-04T08:52:21+00:00"

# 2. Locate the existing "Mapping Table 0" sheet
$ws0 = $wb.Worksheets.Item("Mapping Table 0")

# 3. Create the new "Mapping Table 1" sheet right after "Mapping Table 0"
$ws1 = $wb.Worksheets.Add([Type]::Missing, $ws0)
$ws1.Name = "Mapping Table 1"

# 4. Copy the header row and the group-header row (with their styles) into the new sheet
$ws0.Range("A1:E2").Copy($ws1.Range("A1:E2"))

# 5. Move the last data row (phase-III-IV -> phase-3-phase-4) into the new sheet as row 3
$ws0.Range("A15:E15").Copy($ws1.Range("A3:E3"))

# 6. Fix the group target URI on "Mapping Table 0" (row 2, column D)
$ws0.Cells.Item(2, 4).Value = "http://terminology.hl7.org/CodeSystem/research-study-phase"

# 7. Remove the obsolete phase-IV -> phase-4 row (row 14) from "Mapping Table 0"
$ws0.Rows.Item(14).Delete()

# 8. Remove the row that was just moved to "Mapping Table 1" (now shifted up to row 14)
$ws0.Rows.Item(14).Delete()
